# Update the cached "datetimeFigureOut" date field text, wherever it
# appears (slide master + every slide layout's Date placeholder), from
# 2/20/2019 -> 2/21/19.
$p = $ppt.ActivePresentation
$newDate = "2/21/19"

$m = $p.SlideMaster
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $lay = $m.CustomLayouts.Item($li)
    for ($si = 1; $si -le $lay.Shapes.Count; $si++) {
        $sh = $lay.Shapes.Item($si)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide 1 subtitle: flesh out "Members:" with the actual team roster,
# split into runs the way PowerPoint would after typing + spellcheck.
$s1 = $p.Slides.Item(1)
$sub = $s1.Shapes.Item(2)
$tr = $sub.TextFrame.TextRange
$tr.Text = "Members: Omkar Bhambure, Yueyang Chen, Rachel Gebhart, Isaac Kretzmer"

$tr.Characters(16, 8).Text  = "Bhambure"                 # "Bhambure"
$tr.Characters(24, 2).Text  = ", "                        # ", "
$tr.Characters(26, 7).Text  = "Yueyang"                   # "Yueyang"
$tr.Characters(33, 7).Text  = " Chen, "                   # " Chen, "
$tr.Characters(40, 22).Text = "Rachel Gebhart, Isaac "    # "Rachel Gebhart, Isaac "
$tr.Characters(62, 8).Text  = "Kretzmer"                  # "Kretzmer"
